$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 2
Set-TextValue 'D2' '29.385.92'
Set-TextValue 'E2' '  +0.23%  '

# Row 3
Set-TextValue 'D3' '1.883.91'

# Row 4
Set-TextValue 'E4' '  +0.04%  '

# Row 5
Set-TextValue 'D5' '0.7126'
Set-TextValue 'E5' '  +0.15%  '

# Row 6
Set-TextValue 'D6' '242.35'
Set-TextValue 'E6' '  +0.02%  '

# Row 7
Set-TextValue 'E7' '  +0.04%  '

# Row 8
Set-TextValue 'D8' '0.08033'
Set-TextValue 'E8' '  +3.42%  '

# Row 9
Set-TextValue 'D9' '0.3128'
Set-TextValue 'E9' '  +0.72%  '

# Row 10
Set-TextValue 'D10' '25.26'
Set-TextValue 'E10' '  +1.22%  '

# Row 11
Set-TextValue 'D11' '0.08341'
Set-TextValue 'E11' '  -1.46%  '

# Row 12
Set-TextValue 'D12' '1.894.08'
Set-TextValue 'E12' '  +0.67%  '

# Row 13
Set-TextValue 'B13' 'Polkadot'
Set-TextValue 'C13' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 'D13' '5.248'
Set-TextValue 'E13' '  +0.65%  '

# Row 14
Set-TextValue 'B14' 'Polygon'
Set-TextValue 'C14' 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue 'D14' '0.7207'
Set-TextValue 'E14' '  +1.31%  '

# Row 15
Set-TextValue 'D15' '92.59'
Set-TextValue 'E15' '  +1.28%  '

# Row 16
Set-TextValue 'D16' '6.326'
Set-TextValue 'E16' '  +5.26%  '

# Row 17
Set-TextValue 'D17' '0.000008446'
Set-TextValue 'E17' '  +2.29%  '

# Row 18
Set-TextValue 'D18' '29.407.95'
Set-TextValue 'E18' '  +0.30%  '

# Row 19
Set-TextValue 'D19' '241.14'
Set-TextValue 'E19' '  -0.64%  '

# Row 20
Set-TextValue 'B20' 'Avalanche'
Set-TextValue 'C20' 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue 'D20' '13.26'
Set-TextValue 'E20' '  +0.05%  '

# Row 21
Set-TextValue 'B21' 'WrappedliquidstakedEther2.0'
Set-TextValue 'C21' 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue 'D21' '2.134.18'
Set-TextValue 'E21' '  +0.01%  '

# Row 23
Set-TextValue 'D23' '7.840'
Set-TextValue 'E23' '  -0.01%  '

# Row 24
Set-TextValue 'E24' '  +0.06%  '

# Row 25
Set-TextValue 'D25' '0.1586'
Set-TextValue 'E25' '  -1.84%  '

# Row 26
Set-TextValue 'D26' '163.97'
Set-TextValue 'E26' '  +0.83%  '

# Row 27
Set-TextValue 'D27' '9.055'
Set-TextValue 'E27' '  +0.31%  '

# Row 29
Set-TextValue 'E29' '  -0.30%  '

# Row 30
Set-TextValue 'D30' '4.417'
Set-TextValue 'E30' '  +0.21%  '

# Row 31
Set-TextValue 'D31' '4.344'
Set-TextValue 'E31' '  +0.18%  '

# Row 32
Set-TextValue 'D32' '1.206'
Set-TextValue 'E32' '  -5.63%  '

# Row 33
Set-TextValue 'D33' '0.05370'

# Row 34
Set-TextValue 'E34' '  +0.93%  '

# Row 35
Set-TextValue 'D35' '1.184'
Set-TextValue 'E35' '  +0.48%  '

# Row 36
Set-TextValue 'D36' '0.7493'
Set-TextValue 'E36' '  +1.19%  '

# Row 37
Set-TextValue 'D37' '2.703'
Set-TextValue 'E37' '  +0.55%  '

# Row 38
Set-TextValue 'E38' '  +1.23%  '

# Row 39
Set-TextValue 'D39' '1.287.18'
Set-TextValue 'E39' '  +9.60%  '

# Row 40
Set-TextValue 'E40' '  +0.77%  '

# Row 41
Set-TextValue 'D41' '6.598'
Set-TextValue 'E41' '  +3.22%  '

# Row 42
Set-TextValue 'D42' '0.9258'
Set-TextValue 'E42' '  +4.26%  '

# Row 43
Set-TextValue 'D43' '111.57'
Set-TextValue 'E43' '  +4.76%  '

# Row 44
Set-TextValue 'D44' '73.60'
Set-TextValue 'E44' '  +0.88%  '

# Row 45
Set-TextValue 'E45' '  +0.06%  '

# Row 46
Set-TextValue 'D46' '0.00000000128'
Set-TextValue 'E46' '  +5.78%  '

# Row 47
Set-TextValue 'D47' '2.031.25'
Set-TextValue 'E47' '  +0.10%  '

# Row 48
Set-TextValue 'E48' '  -0.36%  '

# Row 49
Set-TextValue 'D49' '0.5216'
Set-TextValue 'E49' '  +0.25%  '

# Row 50
Set-TextValue 'D50' '9.501'
Set-TextValue 'E50' '  +0.96%  '

# Row 51
Set-TextValue 'D51' '0.4403'
Set-TextValue 'E51' '  +2.07%  '
